$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = "31/12/2014"
$ws.Cells.Item(2, 4).Value = 7.8
$ws.Cells.Item(3, 3).Value = "31/12/2014"
$ws.Cells.Item(3, 4).Value = 8.2
$ws.Cells.Item(4, 3).Value = "31/12/2014"
$ws.Cells.Item(4, 4).Value = 9.2
$ws.Cells.Item(5, 3).Value = "31/12/2014"
$ws.Cells.Item(5, 4).Value = 9.2
$ws.Cells.Item(6, 3).Value = "31/12/2014"
$ws.Cells.Item(6, 4).Value = 9.3
$ws.Cells.Item(7, 3).Value = "31/12/2014"
$ws.Cells.Item(7, 4).Value = 8.9
$ws.Cells.Item(8, 3).Value = "31/12/2014"
$ws.Cells.Item(8, 4).Value = 8.8
$ws.Cells.Item(9, 3).Value = "31/12/2014"
$ws.Cells.Item(9, 4).Value = 7.6
$ws.Cells.Item(10, 3).Value = "31/12/2014"
$ws.Cells.Item(10, 4).Value = 6.3
$ws.Cells.Item(11, 3).Value = "31/12/2014"
$ws.Cells.Item(11, 4).Value = 5.9
$ws.Cells.Item(12, 3).Value = "31/12/2014"
$ws.Cells.Item(12, 4).Value = 5.4
$ws.Cells.Item(13, 3).Value = "31/12/2014"
$ws.Cells.Item(13, 4).Value = 3.6
$ws.Cells.Item(14, 3).Value = "31/12/2014"
$ws.Cells.Item(14, 4).Value = 3.3
$ws.Cells.Item(15, 3).Value = "31/12/2014"
$ws.Cells.Item(15, 4).Value = 2.4
$ws.Cells.Item(16, 3).Value = "31/12/2014"
$ws.Cells.Item(16, 4).Value = 1.6
$ws.Cells.Item(17, 3).Value = "31/12/2014"
$ws.Cells.Item(17, 4).Value = 1.1
$ws.Cells.Item(18, 3).Value = "31/12/2014"
$ws.Cells.Item(18, 4).Value = 1.3
$ws.Cells.Item(19, 3).Value = "31/12/2019"
$ws.Cells.Item(19, 4).Value = 7.5
$ws.Cells.Item(20, 3).Value = "31/12/2019"
$ws.Cells.Item(20, 4).Value = 7.7
$ws.Cells.Item(21, 3).Value = "31/12/2019"
$ws.Cells.Item(21, 4).Value = 8
$ws.Cells.Item(22, 3).Value = "31/12/2019"
$ws.Cells.Item(22, 4).Value = 9
$ws.Cells.Item(23, 3).Value = "31/12/2019"
$ws.Cells.Item(23, 4).Value = 9
$ws.Cells.Item(24, 3).Value = "31/12/2019"
$ws.Cells.Item(24, 4).Value = 8.7
$ws.Cells.Item(25, 3).Value = "31/12/2019"
$ws.Cells.Item(25, 4).Value = 8
$ws.Cells.Item(26, 3).Value = "31/12/2019"
$ws.Cells.Item(26, 4).Value = 8.2
$ws.Cells.Item(27, 3).Value = "31/12/2019"
$ws.Cells.Item(27, 4).Value = 7.5
$ws.Cells.Item(28, 3).Value = "31/12/2019"
$ws.Cells.Item(28, 4).Value = 6.1
$ws.Cells.Item(29, 3).Value = "31/12/2019"
$ws.Cells.Item(29, 4).Value = 5.3
$ws.Cells.Item(30, 3).Value = "31/12/2019"
$ws.Cells.Item(30, 4).Value = 4
$ws.Cells.Item(31, 3).Value = "31/12/2019"
$ws.Cells.Item(31, 4).Value = 3.3
$ws.Cells.Item(32, 3).Value = "31/12/2019"
$ws.Cells.Item(32, 4).Value = 2.6
$ws.Cells.Item(33, 3).Value = "31/12/2019"
$ws.Cells.Item(33, 4).Value = 1.9
$ws.Cells.Item(34, 3).Value = "31/12/2019"
$ws.Cells.Item(34, 4).Value = 1.3
$ws.Cells.Item(35, 3).Value = "31/12/2019"
$ws.Cells.Item(35, 4).Value = 1.6
$ws.Cells.Item(36, 3).Value = "31/12/2023"
$ws.Cells.Item(36, 4).Value = 7.5
$ws.Cells.Item(37, 3).Value = "31/12/2023"
$ws.Cells.Item(37, 4).Value = 8.1
$ws.Cells.Item(38, 3).Value = "31/12/2023"
$ws.Cells.Item(38, 4).Value = 7.1
$ws.Cells.Item(39, 3).Value = "31/12/2023"
$ws.Cells.Item(39, 4).Value = 7.8
$ws.Cells.Item(40, 3).Value = "31/12/2023"
$ws.Cells.Item(40, 4).Value = 8.4
$ws.Cells.Item(41, 3).Value = "31/12/2023"
$ws.Cells.Item(41, 4).Value = 7.9
$ws.Cells.Item(42, 3).Value = "31/12/2023"
$ws.Cells.Item(42, 4).Value = 8.2
$ws.Cells.Item(43, 3).Value = "31/12/2023"
$ws.Cells.Item(43, 4).Value = 7.4
$ws.Cells.Item(44, 3).Value = "31/12/2023"
$ws.Cells.Item(44, 4).Value = 7.3
$ws.Cells.Item(45, 3).Value = "31/12/2023"
$ws.Cells.Item(45, 4).Value = 6.8
$ws.Cells.Item(46, 3).Value = "31/12/2023"
$ws.Cells.Item(46, 4).Value = 5.4
$ws.Cells.Item(47, 3).Value = "31/12/2023"
$ws.Cells.Item(47, 4).Value = 5.3
$ws.Cells.Item(48, 3).Value = "31/12/2023"
$ws.Cells.Item(48, 4).Value = 4.2
$ws.Cells.Item(49, 3).Value = "31/12/2023"
$ws.Cells.Item(49, 4).Value = 2.9
$ws.Cells.Item(50, 3).Value = "31/12/2023"
$ws.Cells.Item(50, 4).Value = 2
$ws.Cells.Item(51, 3).Value = "31/12/2023"
$ws.Cells.Item(51, 4).Value = 1.6
$ws.Cells.Item(52, 3).Value = "31/12/2023"
$ws.Cells.Item(52, 4).Value = 2.1
